# Added a "delete confirm" message scenario to the StaffData sheet:
# append a new first/last name pair as a new row (A2:A3).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StaffData")

$ws.Range("A2").Value = "SFirstLRZST"
$ws.Range("A3").Value = "SLastVSAHI"
